$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("IF function example #1")
$ws2 = $wb.Worksheets.Item("IF function example #2")
$ws3 = $wb.Worksheets.Item("IF function example #3")

# --- Sheet 2 ("IF function example #2"): fill in the "Color" column (R) ---
# R2 header should read "Color" (same label as the B2 header in the example section)
$ws2.Range("R2").Value = "Color"

# R3:R27 mirror the cyclical color pattern already used in column B
# (Blue, Green, Red, White, Black repeating)
$colors = @("Blue", "Green", "Red", "White", "Black")
for ($i = 0; $i -lt 25; $i++) {
    $row = 3 + $i
    $ws2.Cells.Item($row, 18).Value = $colors[$i % 5]
}

# R29 is the exercise cell where the user enters the color to match against
$ws2.Range("R29").Value = "Green"

# --- Selections left behind on each sheet ---
$ws1.Range("B3").Select()
$ws2.Range("R3").Select()
$ws3.Range("S5").Select()

# --- Active sheet ends up being sheet 3 ---
$ws3.Activate()
